$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap all match-data columns (B..AB) between two rows, leaving column A
# (the row index) untouched, since the two fixtures were re-ordered.
function Swap-Rows($row1, $row2) {
    $addr1 = "B" + $row1 + ":AB" + $row1
    $addr2 = "B" + $row2 + ":AB" + $row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-Rows 88 89
Swap-Rows 115 116
